# Apply updated Betfair Back/Lay odds values for 2025-12-29 Jogos do Dia sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.64
$ws.Range("G2").Value = 1.66
$ws.Range("N2").Value = 5
$ws.Range("P2").Value = 2.36
$ws.Range("Q2").Value = 1.68
$ws.Range("R2").Value = 1.54
$ws.Range("T2").Value = 1.75
$ws.Range("U2").Value = 2.22
$ws.Range("W2").Value = 2.5
$ws.Range("X2").Value = 21
$ws.Range("Z2").Value = 50
$ws.Range("AA2").Value = 150

# Row 3
$ws.Range("G3").Value = 2.5
$ws.Range("K3").Value = 3.25
$ws.Range("V3").Value = 1.26

# Row 4
$ws.Range("G4").Value = 2.5
$ws.Range("I4").Value = 3.35
$ws.Range("P4").Value = 2.52
$ws.Range("Q4").Value = 1.52
$ws.Range("R4").Value = 1.56
$ws.Range("V4").Value = 1.43
$ws.Range("W4").Value = 1.66
$ws.Range("Y4").Value = 26
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 15
$ws.Range("AD4").Value = 20
$ws.Range("AF4").Value = 29
$ws.Range("AG4").Value = 18
$ws.Range("AH4").Value = 21
$ws.Range("AN4").Value = 17.5
$ws.Range("AO4").Value = 25

# Row 5
$ws.Range("L5").Value = 1.56
$ws.Range("M5").Value = 1.16
$ws.Range("N5").Value = 2.22
$ws.Range("O5").Value = 1.65
$ws.Range("Q5").Value = 3
$ws.Range("T5").Value = 2.26

# Row 6
$ws.Range("Q6").Value = 1.67

# Row 7
$ws.Range("N7").Value = 1.87
$ws.Range("R7").Value = 1.21

# Row 8
$ws.Range("F8").Value = 1.44
$ws.Range("G8").Value = 1.7
$ws.Range("H8").Value = 6.8
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 3.1
$ws.Range("K8").Value = 6.2
$ws.Range("N8").Value = 2.06
$ws.Range("O8").Value = 1.31
$ws.Range("Q8").Value = 2.28
$ws.Range("S8").Value = 2.96
$ws.Range("V8").Value = 1.01
$ws.Range("W8").Value = 2.44

# Row 9
$ws.Range("H9").Value = 6.8
$ws.Range("I9").Value = 7.2
$ws.Range("S9").Value = 5
$ws.Range("W9").Value = 2.42
$ws.Range("Z9").Value = 55
$ws.Range("AE9").Value = 150
$ws.Range("AG9").Value = 11

# Row 10
$ws.Range("F10").Value = 1.12
$ws.Range("I10").Value = 36
$ws.Range("N10").Value = 6.8
$ws.Range("O10").Value = 1.14
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 1.38
$ws.Range("R10").Value = 1.79
$ws.Range("S10").Value = 1.97
$ws.Range("U10").Value = 1.44
$ws.Range("X10").Value = 980
$ws.Range("Z10").Value = 420
$ws.Range("AD10").Value = 130
$ws.Range("AF10").Value = 8
$ws.Range("AG10").Value = 980
$ws.Range("AI10").Value = 520
$ws.Range("AK10").Value = 980
$ws.Range("AM10").Value = 490
